# Boolean Exempt Process Emissions From Carbon Tax
# Flip the "Exempt Process Emissions from Carbon Tax" control lever from
# TRUE (1) to FALSE (0) on the BEPEfCT sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 0
